$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "tier" column (D): header + formula classifying A by magnitude.
$ws.Range("D1").Value = "tier"
$ws.Range("D2").Formula = "=IF(A2<100,1,IF(A2<10000,2,IF(A2<1000000,3,4)))"
$ws.Range("D3:D32").Formula = "=IF(A3<100,1,IF(A3<10000,2,IF(A3<1000000,3,4)))"

# Move the active selection to D10 (was D1:G1048576 anchored at D1).
$ws.Range("D10").Select() | Out-Null

# Restore the workbook window position recorded at save time.
$win = $wb.Windows.Item(1)
$win.Left = 5520
$win.Top = -20300
